$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-10 (replaces old rows 2-15)
$data = @(
    @("Redes Inalambricas", 49362, "LOPEZ - MUNOZ MAURO ALBERTO", "L", "13:00", "13:59", "1CCO4/301"),
    @("Redes Inalambricas", 49362, "LOPEZ - MUNOZ MAURO ALBERTO", "M", "13:00", "14:59", "1CCO3/303"),
    @("Redes Inalambricas", 49362, "LOPEZ - MUNOZ MAURO ALBERTO", "V", "13:00", "14:59", "1CCO4/301"),
    @("Tec.de Inteligencia Artificial", 49245, "TECUANHUEHUE - VERA PEDRO", "L", "11:00", "11:59", "1CCO3/114"),
    @("Tec.de Inteligencia Artificial", 49245, "TECUANHUEHUE - VERA PEDRO", "M", "11:00", "12:59", "1CCO3/114"),
    @("Tec.de Inteligencia Artificial", 49245, "TECUANHUEHUE - VERA PEDRO", "V", "11:00", "12:59", "1CCO5/202"),
    @("Teoria de Control", 49190, "HERNANDEZ - AMECA JOSE LUIS", "L", "12:00", "12:59", "1CCO3/310"),
    @("Teoria de Control", 49190, "HERNANDEZ - AMECA JOSE LUIS", "A", "11:00", "12:59", "1CCO3/310"),
    @("Teoria de Control", 49190, "HERNANDEZ - AMECA JOSE LUIS", "J", "11:00", "12:59", "1CCO1/002")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
}

# Remove the now-obsolete rows 11-15
$ws.Range("A11:G15").Delete()
